# Applies the LOQ4078.docx reshuffle described by the commit diff.
#
# The paragraph-level formatting (styles, bold/italic run properties)
# stays exactly where it is; only the literal text content moves
# between paragraphs / between runs of one paragraph. We therefore:
#   1. capture every "old" text value we will need, up front, before
#      touching anything (so later writes never clobber a value we
#      still have to read);
#   2. write the "new" values into their destinations.
#
# Word paragraph numbering used below (1-based, confirmed against the
# live object model):
#   6  -> "Este curso visa apresentar..."           (plain run)
#   7  -> "This course aims to introduce..."        (italic run)
#   9  -> "5817066 - Félix Monteiro Pereira"         (ListBullet)
#   11 -> "Serão abordados assuntos..."              (plain run)
#   12 -> "Issues of interest and concern..."        (italic run)
#   14 -> "1. Introdução a modelagem..."             (plain run, has a <w:br/>)
#   17 -> "Método: / Critério: / Norma de recuperação: " bullet paragraph
#   19 -> Bibliography list (6 numbered references)
#
# The full content-move cycle (confirmed against the target XML) is:
#   P6 <- P11 <- P14 <- (P17 "Método:" run)
#                           ^
#                   (P17 "Critério:" run) <- (P17 "Norma de recuperação:" run) <- P19 <- P9 <- P6
#   and, separately: P7 <-> P12 (simple swap)

$d = $word.ActiveDocument

function Chomp([string]$s) {
    # Paragraph.Range.Text always ends with the paragraph mark (CR,
    # char 13); strip it so re-assigning the text doesn't insert an
    # extra paragraph.
    if ($s.Length -gt 0 -and [int][char]$s[$s.Length - 1] -eq 13) {
        return $s.Substring(0, $s.Length - 1)
    }
    return $s
}

# ---------------------------------------------------------------------
# 1. Capture the "before" text of every whole paragraph whose content
#    moves elsewhere.
# ---------------------------------------------------------------------
$p6  = Chomp $d.Paragraphs.Item(6).Range.Text
$p7  = Chomp $d.Paragraphs.Item(7).Range.Text
$p9  = Chomp $d.Paragraphs.Item(9).Range.Text
$p11 = Chomp $d.Paragraphs.Item(11).Range.Text
$p12 = Chomp $d.Paragraphs.Item(12).Range.Text
$p14 = Chomp $d.Paragraphs.Item(14).Range.Text
$p19 = Chomp $d.Paragraphs.Item(19).Range.Text

# ---------------------------------------------------------------------
# 2. Capture the three "content" runs inside paragraph 17 (the bold
#    labels "Método:", "Critério:", "Norma de recuperação:" stay put;
#    the plain-text run that follows each label is what moves).
# ---------------------------------------------------------------------
$p17 = $d.Paragraphs.Item(17)

$rMetodoLabel = $p17.Range.Duplicate
$rMetodoLabel.Find.Execute("M" + [char]0x00E9 + "todo: ") | Out-Null
$afterMetodo = $rMetodoLabel.End

$rCriterioLabel = $p17.Range.Duplicate
$rCriterioLabel.Find.Execute("Crit" + [char]0x00E9 + "rio: ") | Out-Null
$beforeCriterio = $rCriterioLabel.Start
$afterCriterio = $rCriterioLabel.End

$rNormaLabel = $p17.Range.Duplicate
$rNormaLabel.Find.Execute("Norma de recupera" + [char]0x00E7 + [char]0x00E3 + "o: ") | Out-Null
$beforeNorma = $rNormaLabel.Start
$afterNorma = $rNormaLabel.End

$paraEnd = $p17.Range.End

$content1 = $d.Range($afterMetodo, $beforeCriterio).Text     # after "Método:"   up to "Critério:"
$content2 = $d.Range($afterCriterio, $beforeNorma).Text       # after "Critério:" up to "Norma..."
$content3 = $d.Range($afterNorma, $paraEnd).Text              # after "Norma de recuperação:" to para end

# ---------------------------------------------------------------------
# 3. Write the simple whole-paragraph swaps (all sources already cached
#    above, so write order amongst these does not matter).
# ---------------------------------------------------------------------
$d.Paragraphs.Item(6).Range.Text  = $p11
$d.Paragraphs.Item(7).Range.Text  = $p12
$d.Paragraphs.Item(9).Range.Text  = $p6
$d.Paragraphs.Item(11).Range.Text = $p14
$d.Paragraphs.Item(12).Range.Text = $p7
$d.Paragraphs.Item(14).Range.Text = $content1
$d.Paragraphs.Item(19).Range.Text = $p9

# ---------------------------------------------------------------------
# 4. Rewrite paragraph 17's three content runs, back-to-front so the
#    offsets captured in step 2 (which all precede $paraEnd) remain
#    valid while still in use:
#      content1 (first)  <- old content2 (unchanged text, just moved)
#      content2 (middle) <- old content3 + a trailing line break
#      content3 (last)   <- bibliography list (was paragraph 19)
# ---------------------------------------------------------------------
$d.Range($afterNorma, $paraEnd).Text         = $p19
$d.Range($afterCriterio, $beforeNorma).Text  = $content3 + [char]11
$d.Range($afterMetodo, $beforeCriterio).Text = $content2
